$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.089.93'
$ws.Range('E2').Value = '  -4.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.485.61'
$ws.Range('E3').Value = '  -3.67%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.80'
$ws.Range('E5').Value = '  -2.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.56'
$ws.Range('E6').Value = '  -6.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.573'
$ws.Range('E8').Value = '  -3.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.518.94'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0996'
$ws.Range('E10').Value = '  -4.81%  '
$ws.Range('E11').Value = '  -2.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.53'
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  -3.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.928.58'
$ws.Range('E14').Value = '  -3.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.02'
$ws.Range('E15').Value = '  -6.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.996.58'
$ws.Range('E16').Value = '  -4.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.506.06'
$ws.Range('E18').Value = '  -3.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.32'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.29'
$ws.Range('E20').Value = '  -5.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.83'
$ws.Range('E21').Value = '  -4.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.75'
$ws.Range('E23').Value = '  -4.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.25'
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.441'
$ws.Range('E25').Value = '  -10.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.610.05'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.992'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.81'
$ws.Range('E29').Value = '  -4.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.91'
$ws.Range('E30').Value = '  -5.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0776'
$ws.Range('E31').Value = '  -7.59%  '
$ws.Range('E32').Value = '  -6.73%  '
$ws.Range('E33').Value = '  -5.65%  '
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.05'
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.44'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.51'
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.44'
$ws.Range('E38').Value = '  -9.14%  '
$ws.Range('E39').Value = '  -10.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.92'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '308.34'
$ws.Range('E41').Value = '  -6.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.77'
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.67'
$ws.Range('E43').Value = '  -7.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.822'
$ws.Range('E44').Value = '  -10.09%  '
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.597'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.77'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.64'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0930'
$ws.Range('E49').Value = '  -3.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.58'
$ws.Range('E50').Value = '  -5.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0517'
$ws.Range('E51').Value = '  -5.74%  '
